$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4721909463405609
$ws.Range("B1").Value = 1.802385330200195
$ws.Range("C1").Value = 4.878477573394775
$ws.Range("D1").Value = 1.68873393535614
$ws.Range("E1").Value = 0.8877521753311157
